# ---------------------------------------------------------------------------
# Update "template_enrollment.xlsx": replace the old academic-year header row
# and the yellow instructional note with the new enrollment bulk-insert
# layout (Region / HEI Code / Program / Major / Semester / 2017-M / 2017-F).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Get rid of the merged "yellow note" block (I3:L5) -----------------
$ws.Range("I3:L5").UnMerge()
$ws.Rows("3:5").Delete()

# --- 2. Re-label the existing header cells ---------------------------------
$ws.Range("B1").Value = "HEI Code"

# C1/D1/E1 used to hold the bold year numbers (2017/2018/2019). They become
# the new text headers "Program" / "Major" / "Semester", styled the same way
# the "HEI Code" header (B1) already is.
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "Program"
$ws.Range("D1").Value = "Major"
$ws.Range("E1").Value = "Semester"

# F1 used to hold the bold year number 2020; it becomes the (non-bold)
# semester value "2017-M", and a new G1 cell is added for "2017-F" using the
# same formatting.
$ws.Range("F1").Value = "2017-M"
$ws.Range("F1").Font.Bold = $false

$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("G1").Value = "2017-F"

$excel.CutCopyMode = $false

# --- 3. Column widths: columns B through E are now the "wide" columns ------
$ws.Columns("B:E").ColumnWidth = 16.17

# --- 4. Selection lands on D2, matching the refreshed template -------------
$ws.Range("D2").Select() | Out-Null
